# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker table (rows 16-20) is re-ordered: the "ZOBEIDA CASTRO POLO"
# record (previously the 4th data row, row 19) moves up to become the
# first data row (row 16); the other three previously-existing records
# shift down one row, and their "Salario Basico" (column G) amounts are
# updated to new values. The last row (20, DANNES ELENA CASTILLO MORENO)
# is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: ZOBEIDA CASTRO POLO (moved up from the old row 19)
$ws.Range("C16").Value = "52406183"
$ws.Range("D16").Value = "ZOBEIDA CASTRO POLO"
$ws.Range("E16").Value = "2108"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 908528

# New row 17: YANILETH ESTHER MEDINA BETANCOURT (was row 16), new Salario Basico
$ws.Range("C17").Value = "1047434072"
$ws.Range("D17").Value = "YANILETH ESTHER MEDINA BETANCOURT"
$ws.Range("E17").Value = "2202"
$ws.Range("F17").Value = 20000
$ws.Range("G17").Value = 1500000

# New row 18: CARLOS ENRIQUE QUINTANA SLAGADO (was row 17), new Salario Basico
$ws.Range("C18").Value = "1082241049"
$ws.Range("D18").Value = "CARLOS ENRIQUE QUINTANA SLAGADO"
$ws.Range("E18").Value = "2202"
$ws.Range("F18").Value = 18666
$ws.Range("G18").Value = 2104000

# New row 19: JOSE DANIEL VARGAS PAJARO (was row 18), new Salario Basico
$ws.Range("C19").Value = "1007739316"
$ws.Range("D19").Value = "JOSE DANIEL VARGAS PAJARO"
$ws.Range("E19").Value = "2202"
$ws.Range("F19").Value = 20000
$ws.Range("G19").Value = 1600000

# Row 20: DANNES ELENA CASTILLO MORENO - unchanged, left as-is
